# Update "想去人数" (number of people interested) figures on the 展览
# (Exhibition) and 全部类型 (All types) sheets, mirroring the scraped data
# refresh captured in the gh-pages output regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 14502
$wsExhibit.Range("F3").Value = 340
$wsExhibit.Range("F4").Value = 696
$wsExhibit.Range("F5").Value = 242
$wsExhibit.Range("F6").Value = 584
$wsExhibit.Range("F7").Value = 1531
$wsExhibit.Range("F8").Value = 145

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 14502
$wsAll.Range("F3").Value = 340
$wsAll.Range("F4").Value = 696
$wsAll.Range("F5").Value = 242
$wsAll.Range("F8").Value = 584
$wsAll.Range("F9").Value = 1531
$wsAll.Range("F11").Value = 145
